$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update D3: change "Dropdown" wording to "Select" for the ticket/frequency attribute definition.
$ws.Range("D3").Value = "Ticketart:Select(Bus,Zug,U-Bahn);Häufigkeit:Select(Täglich,Wöchentlich,Selten)"

# Move active selection to D3.
$ws.Range("D3").Select()
